$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data as scraped on Sat May  6 10:52:55 UTC 2023
# A leading apostrophe forces Excel to store the value as literal text
# (matches the workbooks inlineStr cells) instead of auto-coercing
# numeric-looking strings ("0.07060", "7.640", "325.19", ...) into numbers.

$ws.Range('D2').Value = "'29.267.53"
$ws.Range('E2').Value = "'  +0.33%  "

$ws.Range('D3').Value = "'1.931.30"
$ws.Range('E3').Value = "'  +1.30%  "

$ws.Range('E4').Value = "'  -0.04%  "

$ws.Range('D5').Value = "'325.19"
$ws.Range('E5').Value = "'  -0.26%  "

$ws.Range('D6').Value = "'0.9996"
$ws.Range('E6').Value = "'  -0.19%  "

$ws.Range('D7').Value = "'0.4618"
$ws.Range('E7').Value = "'  +0.10%  "

$ws.Range('D8').Value = "'0.3866"
$ws.Range('E8').Value = "'  -0.70%  "

$ws.Range('B9').Value = "'OKB"
$ws.Range('C9').Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('D9').Value = "'45.79"
$ws.Range('E9').Value = "'  -0.89%  "

$ws.Range('B10').Value = "'Dogecoin"
$ws.Range('C10').Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range('D10').Value = "'0.07795"
$ws.Range('E10').Value = "'  -1.14%  "

$ws.Range('B11').Value = "'Polygon"
$ws.Range('C11').Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range('D11').Value = "'0.9706"
$ws.Range('E11').Value = "'  -2.04%  "

$ws.Range('B12').Value = "'Solana"
$ws.Range('C12').Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range('D12').Value = "'22.56"
$ws.Range('E12').Value = "'  +2.52%  "

$ws.Range('B13').Value = "'WrappedEther"
$ws.Range('C13').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D13').Value = "'1.920.79"
$ws.Range('E13').Value = "'  +0.57%  "

$ws.Range('B14').Value = "'Chainlink"
$ws.Range('C14').Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range('D14').Value = "'7.066"
$ws.Range('E14').Value = "'  +0.26%  "

$ws.Range('B15').Value = "'Polkadot"
$ws.Range('C15').Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('D15').Value = "'5.768"
$ws.Range('E15').Value = "'  -0.04%  "

$ws.Range('B16').Value = "'TRON"
$ws.Range('C16').Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range('D16').Value = "'0.07060"
$ws.Range('E16').Value = "'  +0.31%  "

$ws.Range('B17').Value = "'Litecoin"
$ws.Range('C17').Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range('D17').Value = "'86.67"
$ws.Range('E17').Value = "'  -1.67%  "

$ws.Range('B18').Value = "'BinanceUSD"
$ws.Range('C18').Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range('D18').Value = "'1.003"
$ws.Range('E18').Value = "'  -0.05%  "

$ws.Range('B19').Value = "'ShibaInu"
$ws.Range('C19').Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range('D19').Value = "'0.000009600"
$ws.Range('E19').Value = "'  -3.56%  "

$ws.Range('B20').Value = "'Avalanche"
$ws.Range('C20').Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('D20').Value = "'17.01"
$ws.Range('E20').Value = "'  -0.45%  "

$ws.Range('B21').Value = "'Dai"
$ws.Range('C21').Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range('D21').Value = "'1.001"
$ws.Range('E21').Value = "'  -0.11%  "

$ws.Range('B22').Value = "'WrappedBTC"
$ws.Range('C22').Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range('D22').Value = "'29.296.31"
$ws.Range('E22').Value = "'  +0.36%  "

$ws.Range('B23').Value = "'Uniswap"
$ws.Range('C23').Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range('D23').Value = "'5.467"
$ws.Range('E23').Value = "'  +2.68%  "

$ws.Range('B24').Value = "'Cosmos"
$ws.Range('C24').Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range('D24').Value = "'11.04"
$ws.Range('E24').Value = "'  -0.88%  "

$ws.Range('B25').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('C25').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D25').Value = "'2.163.22"
$ws.Range('E25').Value = "'  +1.24%  "

$ws.Range('B26').Value = "'Toncoin"
$ws.Range('C26').Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('D26').Value = "'2.091"
$ws.Range('E26').Value = "'  -0.54%  "

$ws.Range('B27').Value = "'Monero"
$ws.Range('C27').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('D27').Value = "'156.89"
$ws.Range('E27').Value = "'  +0.41%  "

$ws.Range('B28').Value = "'EthereumClassic"
$ws.Range('C28').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('D28').Value = "'19.37"
$ws.Range('E28').Value = "'  -0.18%  "

$ws.Range('B29').Value = "'InternetComputer(DFINITY)"
$ws.Range('C29').Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range('D29').Value = "'5.733"
$ws.Range('E29').Value = "'  -3.11%  "

$ws.Range('B30').Value = "'BitcoinCash"
$ws.Range('C30').Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range('D30').Value = "'118.38"
$ws.Range('E30').Value = "'  -0.41%  "

$ws.Range('B31').Value = "'LidoDAOToken"
$ws.Range('C31').Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range('D31').Value = "'1.849"
$ws.Range('E31').Value = "'  -1.64%  "

$ws.Range('B32').Value = "'Stellar"
$ws.Range('C32').Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range('D32').Value = "'0.09338"
$ws.Range('E32').Value = "'  -0.19%  "

$ws.Range('B33').Value = "'ImmutableX"
$ws.Range('C33').Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('D33').Value = "'0.8615"
$ws.Range('E33').Value = "'  -3.83%  "

$ws.Range('B34').Value = "'Filecoin"
$ws.Range('C34').Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('D34').Value = "'5.159"
$ws.Range('E34').Value = "'  -1.36%  "

$ws.Range('B35').Value = "'ARBITRUM"
$ws.Range('C35').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('D35').Value = "'1.303"
$ws.Range('E35').Value = "'  -1.54%  "

$ws.Range('B36').Value = "'HuobiToken"
$ws.Range('C36').Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D36').Value = "'3.074"
$ws.Range('E36').Value = "'  -2.69%  "

$ws.Range('B37').Value = "'Hedera"
$ws.Range('C37').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('D37').Value = "'0.05762"
$ws.Range('E37').Value = "'  -0.61%  "

$ws.Range('B38').Value = "'TrustWalletToken"
$ws.Range('C38').Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range('D38').Value = "'1.153"
$ws.Range('E38').Value = "'  -1.81%  "

$ws.Range('B39').Value = "'VeChain"
$ws.Range('C39').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('D39').Value = "'0.02075"
$ws.Range('E39').Value = "'  -0.62%  "

$ws.Range('B40').Value = "'FraxShare"
$ws.Range('C40').Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('D40').Value = "'7.640"
$ws.Range('E40').Value = "'  -0.65%  "

$ws.Range('B41').Value = "'TheSandbox"
$ws.Range('C41').Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range('D41').Value = "'0.5642"
$ws.Range('E41').Value = "'  -1.27%  "

$ws.Range('B42').Value = "'PEPE"
$ws.Range('C42').Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range('D42').Value = "'0.000003095"
$ws.Range('E42').Value = "'  +55.82%  "

$ws.Range('B43').Value = "'Algorand"
$ws.Range('C43').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('D43').Value = "'0.1775"
$ws.Range('E43').Value = "'  -1.81%  "

$ws.Range('B44').Value = "'Aptos"
$ws.Range('C44').Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range('D44').Value = "'9.340"
$ws.Range('E44').Value = "'  -3.95%  "

$ws.Range('B45').Value = "'MXToken"
$ws.Range('C45').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D45').Value = "'2.713"
$ws.Range('E45').Value = "'  +6.41%  "

$ws.Range('B46').Value = "'Decentraland"
$ws.Range('C46').Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range('D46').Value = "'0.5262"
$ws.Range('E46').Value = "'  -1.88%  "

$ws.Range('B47').Value = "'EnergySwap"
$ws.Range('C47').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('D47').Value = "'11.53"
$ws.Range('E47').Value = "'  -3.26%  "

$ws.Range('B48').Value = "'Cronos"
$ws.Range('C48').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('D48').Value = "'0.06860"
$ws.Range('E48').Value = "'  -2.24%  "

$ws.Range('B49').Value = "'RenderToken"
$ws.Range('C49').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('D49').Value = "'2.074"
$ws.Range('E49').Value = "'  -4.56%  "

$ws.Range('B50').Value = "'NEARProtocol"
$ws.Range('C50').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('D50').Value = "'1.808"
$ws.Range('E50').Value = "'  -1.92%  "

$ws.Range('B51').Value = "'Quant"
$ws.Range('C51').Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range('D51').Value = "'111.09"
$ws.Range('E51').Value = "'  -1.87%  "

